$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 1983.5536
$ws_ALC.Range("J17").Value = 1739.8704
$ws_ALC.Range("L17").Value = 5219.6112
$ws_ALC.Range("N17").Value = -5555.6112
$ws_ALC.Range("H76").Value = 3249.5
$ws_ALC.Range("I76").Value = 2999
$ws_ALC.Range("J76").Value = 3500
$ws_ALC.Range("K76").Value = 2999
$ws_ALC.Range("L76").Value = 3500
$ws_ALC.Range("M76").Value = -2684
$ws_ALC.Range("N76").Value = -4130
$ws_ALC.Range("H79").Value = 3249.5
$ws_ALC.Range("I79").Value = 2999
$ws_ALC.Range("J79").Value = 3500
$ws_ALC.Range("K79").Value = 2999
$ws_ALC.Range("L79").Value = 3500
$ws_ALC.Range("M79").Value = -1907
$ws_ALC.Range("N79").Value = -5684
$ws_ALC.Range("H98").Value = 2901.875
$ws_ALC.Range("I98").Value = 2650
$ws_ALC.Range("J98").Value = 4665
$ws_ALC.Range("K98").Value = 2650
$ws_ALC.Range("L98").Value = 4665
$ws_ALC.Range("M98").Value = -1152
$ws_ALC.Range("N98").Value = -7661
$ws_ALC.Range("H122").Value = 2901.875
$ws_ALC.Range("I122").Value = 2650
$ws_ALC.Range("J122").Value = 4665
$ws_ALC.Range("K122").Value = 7950
$ws_ALC.Range("L122").Value = 13995
$ws_ALC.Range("M122").Value = -5500
$ws_ALC.Range("N122").Value = -18895
$ws_ALC.Range("H137").Value = 54319.105
$ws_ALC.Range("I137").Value = 1521.5714
$ws_ALC.Range("J137").Value = 85117.664
$ws_ALC.Range("K137").Value = 4564.7142
$ws_ALC.Range("L137").Value = 255352.992
$ws_ALC.Range("M137").Value = -2014.7142
$ws_ALC.Range("N137").Value = -260452.992
$ws_ALC.Range("H138").Value = 1792.6207
$ws_ALC.Range("I138").Value = 1489.9517
$ws_ALC.Range("K138").Value = 4469.855100000001
$ws_ALC.Range("M138").Value = 670.1448999999993
$ws_ALC.Range("H141").Value = 801734.4399999999
$ws_ALC.Range("I141").Value = 875851.5600000001
$ws_ALC.Range("K141").Value = 2627554.68
$ws_ALC.Range("M141").Value = -2622374.68

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 5661.3335
$ws_ARM.Range("I32").Value = 5204.0264
$ws_ARM.Range("J32").Value = 8143.857
$ws_ARM.Range("K32").Value = 5204.0264
$ws_ARM.Range("L32").Value = 8143.857
$ws_ARM.Range("M32").Value = -4917.0264
$ws_ARM.Range("N32").Value = -8717.857
$ws_ARM.Range("H45").Value = 1420.5714
$ws_ARM.Range("I45").Value = 1094.5555
$ws_ARM.Range("K45").Value = 1094.5555
$ws_ARM.Range("M45").Value = -717.5554999999999
$ws_ARM.Range("H61").Value = 5926.3125
$ws_ARM.Range("J61").Value = 16998.5
$ws_ARM.Range("L61").Value = 16998.5
$ws_ARM.Range("N61").Value = -17422.5
$ws_ARM.Range("H74").Value = 1730.1765
$ws_ARM.Range("I74").Value = 1234.6666
$ws_ARM.Range("J74").Value = 2919.4
$ws_ARM.Range("K74").Value = 1234.6666
$ws_ARM.Range("L74").Value = 2919.4
$ws_ARM.Range("M74").Value = -360.6666
$ws_ARM.Range("N74").Value = -4667.4
$ws_ARM.Range("H77").Value = 1730.1765
$ws_ARM.Range("I77").Value = 1234.6666
$ws_ARM.Range("J77").Value = 2919.4
$ws_ARM.Range("K77").Value = 6173.333000000001
$ws_ARM.Range("L77").Value = 14597
$ws_ARM.Range("M77").Value = -1805.333000000001
$ws_ARM.Range("N77").Value = -23333
$ws_ARM.Range("H96").Value = 0
$ws_ARM.Range("J96").Value = 0
$ws_ARM.Range("L96").Value = 0
$ws_ARM.Range("N96").ClearContents()
$ws_ARM.Range("H110").Value = 195.2
$ws_ARM.Range("I110").Value = 120.1
$ws_ARM.Range("K110").Value = 120.1
$ws_ARM.Range("M110").Value = 1924.9
$ws_ARM.Range("H132").Value = 1820.2653
$ws_ARM.Range("I132").Value = 1358.475
$ws_ARM.Range("K132").Value = 4075.425
$ws_ARM.Range("M132").Value = -1545.425
$ws_ARM.Range("H136").Value = 5926.3125
$ws_ARM.Range("J136").Value = 16998.5
$ws_ARM.Range("L136").Value = 50995.5
$ws_ARM.Range("N136").Value = -56095.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H80").Value = 8165.5386
$ws_BSM.Range("I80").Value = 299.2
$ws_BSM.Range("J80").Value = 13082
$ws_BSM.Range("K80").Value = 299.2
$ws_BSM.Range("L80").Value = 13082
$ws_BSM.Range("M80").Value = 698.8
$ws_BSM.Range("N80").Value = -15078
$ws_BSM.Range("H83").Value = 8165.5386
$ws_BSM.Range("I83").Value = 299.2
$ws_BSM.Range("J83").Value = 13082
$ws_BSM.Range("K83").Value = 1496
$ws_BSM.Range("L83").Value = 65410
$ws_BSM.Range("M83").Value = 3496
$ws_BSM.Range("N83").Value = -75394
$ws_BSM.Range("H86").Value = 1094363.1
$ws_BSM.Range("I86").Value = 1253749.2
$ws_BSM.Range("K86").Value = 1253749.2
$ws_BSM.Range("M86").Value = -1252626.2
$ws_BSM.Range("H89").Value = 1094363.1
$ws_BSM.Range("I89").Value = 1253749.2
$ws_BSM.Range("K89").Value = 6268746
$ws_BSM.Range("M89").Value = -6263130
$ws_BSM.Range("H134").Value = 4304.7837
$ws_BSM.Range("I134").Value = 4382.9688
$ws_BSM.Range("J134").Value = 3804.4
$ws_BSM.Range("K134").Value = 13148.9064
$ws_BSM.Range("L134").Value = 11413.2
$ws_BSM.Range("M134").Value = -10613.9064
$ws_BSM.Range("N134").Value = -16483.2

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 1800.4595
$ws_CRP.Range("I31").Value = 1290.0714
$ws_CRP.Range("J31").Value = 3388.3333
$ws_CRP.Range("K31").Value = 1290.0714
$ws_CRP.Range("L31").Value = 3388.3333
$ws_CRP.Range("M31").Value = -995.0714
$ws_CRP.Range("N31").Value = -3978.3333
$ws_CRP.Range("H34").Value = 1800.4595
$ws_CRP.Range("I34").Value = 1290.0714
$ws_CRP.Range("J34").Value = 3388.3333
$ws_CRP.Range("K34").Value = 1290.0714
$ws_CRP.Range("L34").Value = 3388.3333
$ws_CRP.Range("M34").Value = -1088.0714
$ws_CRP.Range("N34").Value = -3792.3333
$ws_CRP.Range("H56").Value = 14280
$ws_CRP.Range("I56").Value = 0
$ws_CRP.Range("J56").Value = 14280
$ws_CRP.Range("K56").Value = 0
$ws_CRP.Range("L56").Value = 14280
$ws_CRP.Range("M56").ClearContents()
$ws_CRP.Range("N56").Value = -15970
$ws_CRP.Range("H58").Value = 967198.6
$ws_CRP.Range("I58").Value = 1279487.4
$ws_CRP.Range("J58").Value = 1942.5454
$ws_CRP.Range("K58").Value = 1279487.4
$ws_CRP.Range("L58").Value = 1942.5454
$ws_CRP.Range("M58").Value = -1279284.4
$ws_CRP.Range("N58").Value = -2348.5454
$ws_CRP.Range("H132").Value = 1758.9111
$ws_CRP.Range("I132").Value = 1251.9
$ws_CRP.Range("J132").Value = 2772.9333
$ws_CRP.Range("K132").Value = 3755.7
$ws_CRP.Range("L132").Value = 8318.7999
$ws_CRP.Range("M132").Value = -1225.7
$ws_CRP.Range("N132").Value = -13378.7999
$ws_CRP.Range("H134").Value = 1666.2307
$ws_CRP.Range("I134").Value = 1609.8379
$ws_CRP.Range("J134").Value = 1805.3334
$ws_CRP.Range("K134").Value = 4829.5137
$ws_CRP.Range("L134").Value = 5416.0002
$ws_CRP.Range("M134").Value = -2294.5137
$ws_CRP.Range("N134").Value = -10486.0002
$ws_CRP.Range("H136").Value = 967198.6
$ws_CRP.Range("I136").Value = 1279487.4
$ws_CRP.Range("J136").Value = 1942.5454
$ws_CRP.Range("K136").Value = 3838462.2
$ws_CRP.Range("L136").Value = 5827.6362
$ws_CRP.Range("M136").Value = -3835912.2
$ws_CRP.Range("N136").Value = -10927.6362

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H37").Value = 64500
$ws_CUL.Range("J37").Value = 64500
$ws_CUL.Range("L37").Value = 193500
$ws_CUL.Range("N37").Value = -193724
$ws_CUL.Range("H63").Value = 26035.5
$ws_CUL.Range("I63").Value = 21
$ws_CUL.Range("J63").Value = 52050
$ws_CUL.Range("K63").Value = 63
$ws_CUL.Range("L63").Value = 156150
$ws_CUL.Range("M63").Value = 686
$ws_CUL.Range("N63").Value = -157648
$ws_CUL.Range("H66").Value = 26035.5
$ws_CUL.Range("I66").Value = 21
$ws_CUL.Range("J66").Value = 52050
$ws_CUL.Range("K66").Value = 189
$ws_CUL.Range("L66").Value = 468450
$ws_CUL.Range("M66").Value = 3555
$ws_CUL.Range("N66").Value = -475938
$ws_CUL.Range("H131").Value = 21774564
$ws_CUL.Range("J131").Value = 47737.47
$ws_CUL.Range("L131").Value = 143212.41
$ws_CUL.Range("N131").Value = -153292.41
$ws_CUL.Range("H132").Value = 1300
$ws_CUL.Range("J132").Value = 1300
$ws_CUL.Range("L132").Value = 11700
$ws_CUL.Range("N132").Value = -16760

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H26").Value = 30014
$ws_GSM.Range("J26").Value = 30014
$ws_GSM.Range("L26").Value = 30014
$ws_GSM.Range("N26").Value = -30574
$ws_GSM.Range("H50").Value = 30014
$ws_GSM.Range("J50").Value = 30014
$ws_GSM.Range("L50").Value = 30014
$ws_GSM.Range("N50").Value = -31010
$ws_GSM.Range("H122").Value = 1805.7693
$ws_GSM.Range("I122").Value = 1770.6
$ws_GSM.Range("J122").Value = 1923
$ws_GSM.Range("K122").Value = 5311.799999999999
$ws_GSM.Range("L122").Value = 5769
$ws_GSM.Range("M122").Value = -2861.799999999999
$ws_GSM.Range("N122").Value = -10669
$ws_GSM.Range("H140").Value = 48386.273
$ws_GSM.Range("J140").Value = 48386.273
$ws_GSM.Range("L140").Value = 48386.273
$ws_GSM.Range("N140").Value = -58746.273

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 12214.308
$ws_LTW.Range("I40").Value = 15374
$ws_LTW.Range("K40").Value = 15374
$ws_LTW.Range("M40").Value = -15238
$ws_LTW.Range("H46").Value = 1708.0625
$ws_LTW.Range("I46").Value = 855.9
$ws_LTW.Range("K46").Value = 855.9
$ws_LTW.Range("M46").Value = -667.9

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H62").Value = 6333.1113
$ws_WVR.Range("I62").Value = 6283.1665
$ws_WVR.Range("K62").Value = 6283.1665
$ws_WVR.Range("M62").Value = -5659.1665
$ws_WVR.Range("H65").Value = 6333.1113
$ws_WVR.Range("I65").Value = 6283.1665
$ws_WVR.Range("K65").Value = 31415.8325
$ws_WVR.Range("M65").Value = -28295.8325
$ws_WVR.Range("H122").Value = 45857.715
$ws_WVR.Range("I122").Value = 63641.24
$ws_WVR.Range("J122").Value = 1398.9
$ws_WVR.Range("K122").Value = 190923.72
$ws_WVR.Range("L122").Value = 4196.700000000001
$ws_WVR.Range("M122").Value = -188473.72
$ws_WVR.Range("N122").Value = -9096.700000000001
$ws_WVR.Range("H132").Value = 1640.4546
$ws_WVR.Range("I132").Value = 1105.3235
$ws_WVR.Range("K132").Value = 3315.9705
$ws_WVR.Range("M132").Value = -785.9704999999999
$ws_WVR.Range("H136").Value = 11576104
$ws_WVR.Range("J136").Value = 2323.75
$ws_WVR.Range("L136").Value = 6971.25
$ws_WVR.Range("N136").Value = -12071.25
